$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📙" = "+3"
    "📗" = "✅"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$needsText = @{
    "-3" = $true
    "+3" = $true
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $newVal = $map[$val]
        if ($needsText.ContainsKey($newVal)) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $newVal
    }
}
